$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Raspberry Pi zero" block (columns K:P) ---
# Order of entry matters for shared-string table layout, so we follow the
# same sequence the author appears to have typed things in.

$ws.Range("O7").Value = "prog (lots of tapping)"
$ws.Range("K1").Value = "Raspberry Pi zero"
$ws.Range("O1").Value = "Idel or Program running?"
$ws.Range("P1").Value = "Realtime"
$ws.Range("P7").Value = "buffer overload"

$ws.Range("K4").Value = "Wh"
$ws.Range("L4").Value = "m"
$ws.Range("M4").Value = "s"

$ws.Range("K5").Value = 0.05
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 3
$ws.Range("N5").Formula = "=(K5*3600)/(L5*60 +M5)"
$ws.Range("O5").Value = "idle"

$ws.Range("K6").Value = 0.1438
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = 45
$ws.Range("N6").Formula = "=(K6*3600)/(L6*60 +M6)"
$ws.Range("O6").Value = "prog (only detector)"

$ws.Range("K7").Value = 0.069
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = 0
$ws.Range("N7").Formula = "=(K7*3600)/(L7*60 +M7)"

# --- Re-enter the E14:E21 formulas as one pass, so Excel collapses them into a shared formula group ---
$ws.Range("E14:E21").Formula = "=(B14*3600)/(C14*60 +D14)"

# --- Column widths to fit the new "Raspberry Pi zero" table (matches what Excel's
#     own best-fit auto-sizing produced for the author) ---
$ws.Columns("K:K").ColumnWidth = 14.166666666666666
$ws.Columns("N:N").ColumnWidth = 11.333333333333334
$ws.Columns("O:O").ColumnWidth = 19.666666666666668
$ws.Columns("P:P").ColumnWidth = 12.333333333333334

# --- Selection change ---
$ws.Range("R9").Select()
